$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 10.76772829162301

# Row 3
$ws.Range("B3").Value = 0.04271373187048222
$ws.Range("C3").Value = 0.04071648406533734
$ws.Range("D3").Value = 0.1494219747398047
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 10.42530519761218

# Row 4
$ws.Range("B4").Value = 1.455362044514542
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.1494219747398047
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 3.754798637575387
